$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in the sheet.
$lastRow = $ws.UsedRange.Rows.Count

# Column C holds the "Förändrad" (changed) date; bump every data row's
# date serial from 46061 (2026-02-08) to 46062 (2026-02-09).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    if ($current -eq 46061) {
        $cell.Value2 = 46062
    }
}
